$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.156963590373884
$ws.Range("D2").Value = 4.009826208508815
$ws.Range("E2").Value = 13.0517797733534
$ws.Range("F2").Value = 19.88463245665361
$ws.Range("G2").Value = 21.30854860154199
$ws.Range("H2").Value = 12.30919095588805
$ws.Range("I2").Value = 18.63605908659317
$ws.Range("K2").Value = 10.41413169974733
$ws.Range("M2").Value = 12.99157281482328
$ws.Range("N2").Value = 17.43401136142691
$ws.Range("O2").Value = 17.74895436310258

$ws.Range("B3").Value = 6.039166275242645
$ws.Range("D3").Value = 3.928789575892174
$ws.Range("E3").Value = 12.90363368862676
$ws.Range("F3").Value = 19.85760895936672
$ws.Range("G3").Value = 21.23736206715618
$ws.Range("H3").Value = 12.34233629294067
$ws.Range("I3").Value = 18.72918760439642
$ws.Range("K3").Value = 10.07047087758454
$ws.Range("M3").Value = 12.74849099505098
$ws.Range("N3").Value = 17.48296354965862
$ws.Range("O3").Value = 17.78365092198452

$ws.Range("B4").Value = 5.966593858109411
$ws.Range("D4").Value = 3.8774927197596
$ws.Range("E4").Value = 12.81663994264883
$ws.Range("F4").Value = 19.84680445586836
$ws.Range("G4").Value = 21.20165799374684
$ws.Range("H4").Value = 12.36474362255763
$ws.Range("I4").Value = 18.78931955216837
$ws.Range("K4").Value = 9.851536957735306
$ws.Range("M4").Value = 12.59973984192309
$ws.Range("N4").Value = 17.51462047616478
$ws.Range("O4").Value = 17.80917967898423

$ws.Range("B5").Value = 5.936998827666138
$ws.Range("D5").Value = 3.856216597092263
$ws.Range("E5").Value = 12.78223047599266
$ws.Range("F5").Value = 19.84386088018305
$ws.Range("G5").Value = 21.18913253793329
$ws.Range("H5").Value = 12.3743914698523
$ws.Range("I5").Value = 18.81456767891138
$ws.Range("K5").Value = 9.760411047132564
$ws.Range("M5").Value = 12.53933344100266
$ws.Range("N5").Value = 17.52792425576964
$ws.Range("O5").Value = 17.82064303846904

$ws.Range("B6").Value = 5.932084510643879
$ws.Range("D6").Value = 3.852661685809312
$ws.Range("E6").Value = 12.77658090866513
$ws.Range("F6").Value = 19.84346032883706
$ws.Range("G6").Value = 21.18717522549473
$ws.Range("H6").Value = 12.37602468415592
$ws.Range("I6").Value = 18.81880509507711
$ws.Range("K6").Value = 9.745167018984068
$ws.Range("M6").Value = 12.52931823476352
$ws.Range("N6").Value = 17.53015772677356
$ws.Range("O6").Value = 17.822610482693

$ws.Range("B7").Value = 5.966194761057353
$ws.Range("D7").Value = 3.877207269270719
$ws.Range("E7").Value = 12.81617161357677
$ws.Range("F7").Value = 19.84675884484666
$ws.Range("G7").Value = 21.20148086230888
$ws.Range("H7").Value = 12.3648716452206
$ws.Range("I7").Value = 18.78965704279931
$ws.Range("K7").Value = 9.850315611872842
$ws.Range("M7").Value = 12.59892421358052
$ws.Range("N7").Value = 17.51479826124876
$ws.Range("O7").Value = 17.80932998836585

$ws.Range("B8").Value = 6.116421659657455
$ws.Range("D8").Value = 3.98221316277304
$ws.Range("E8").Value = 12.99990206630625
$ws.Range("F8").Value = 19.87411621117602
$ws.Range("G8").Value = 21.28234989248278
$ws.Range("H8").Value = 12.32019257561326
$ws.Range("I8").Value = 18.6675585381348
$ws.Range("K8").Value = 10.29733006489887
$ws.Range("M8").Value = 12.90770443192352
$ws.Range("N8").Value = 17.45055856268764
$ws.Range("O8").Value = 17.76003961028285

$ws.Range("B9").Value = 6.407361835110563
$ws.Range("D9").Value = 4.175293786994875
$ws.Range("E9").Value = 13.38966581705638
$ws.Range("F9").Value = 19.97345445377833
$ws.Range("G9").Value = 21.50379677539466
$ws.Range("H9").Value = 12.24890339062063
$ws.Range("I9").Value = 18.45144826498498
$ws.Range("K9").Value = 11.10777989615539
$ws.Range("M9").Value = 13.5135565955057
$ws.Range("N9").Value = 17.33723626222366
$ws.Range("O9").Value = 17.69699140797045

$ws.Range("B10").Value = 6.616685792972475
$ws.Range("D10").Value = 4.308562193084458
$ws.Range("E10").Value = 13.69115148915974
$ws.Range("F10").Value = 20.0738853321676
$ws.Range("G10").Value = 21.70373188836988
$ws.Range("H10").Value = 12.20649734960409
$ws.Range("I10").Value = 18.3067674118731
$ws.Range("K10").Value = 11.65909153874004
$ws.Range("M10").Value = 13.95414308673544
$ws.Range("N10").Value = 17.26162878879986
$ws.Range("O10").Value = 17.67126115027506

$ws.Range("B11").Value = 6.710516914360356
$ws.Range("D11").Value = 4.367180480299936
$ws.Range("E11").Value = 13.83096058294531
$ws.Range("F11").Value = 20.12541728742267
$ws.Range("G11").Value = 21.80248015814496
$ws.Range("H11").Value = 12.18937381731668
$ws.Range("I11").Value = 18.24398271421697
$ws.Range("K11").Value = 11.89966030959876
$ws.Range("M11").Value = 14.15264627319414
$ws.Range("N11").Value = 17.2288810984485
$ws.Range("O11").Value = 17.66404185059977

$ws.Range("B12").Value = 6.745814572976845
$ws.Range("D12").Value = 4.389078260163137
$ws.Range("E12").Value = 13.88423202791027
$ws.Range("F12").Value = 20.14575961270876
$ws.Range("G12").Value = 21.84096591711011
$ws.Range("H12").Value = 12.18320143771973
$ws.Range("I12").Value = 18.22064174751279
$ws.Range("K12").Value = 11.98924095354251
$ws.Range("M12").Value = 14.22746372668383
$ws.Range("N12").Value = 17.21671621213302
$ws.Range("O12").Value = 17.66195368336478

$ws.Range("B13").Value = 6.738223525624676
$ws.Range("D13").Value = 4.384375679883409
$ws.Range("E13").Value = 13.87274530514301
$ws.Range("F13").Value = 20.14134190029332
$ws.Range("G13").Value = 21.83262924849463
$ws.Range("H13").Value = 12.18451689266007
$ws.Range("I13").Value = 18.22564935592083
$ws.Range("K13").Value = 11.97001633833175
$ws.Range("M13").Value = 14.21136723030924
$ws.Range("N13").Value = 17.21932565679254
$ws.Range("O13").Value = 17.66237468752063

$ws.Range("B14").Value = 6.713425754706178
$ws.Range("D14").Value = 4.368988091321977
$ws.Range("E14").Value = 13.8353369765202
$ws.Range("F14").Value = 20.12707433767707
$ws.Range("G14").Value = 21.80562469713985
$ws.Range("H14").Value = 12.18885975884203
$ws.Range("I14").Value = 18.24205374559281
$ws.Range("K14").Value = 11.90706082522127
$ws.Range("M14").Value = 14.15880900830249
$ws.Range("N14").Value = 17.22787556169046
$ws.Range("O14").Value = 17.66385711439282

$ws.Range("B15").Value = 6.698204933803247
$ws.Range("D15").Value = 4.359523401911447
$ws.Range("E15").Value = 13.81246451046045
$ws.Range("F15").Value = 20.11844252425922
$ws.Range("G15").Value = 21.78922492078095
$ws.Range("H15").Value = 12.191560520357
$ws.Range("I15").Value = 18.25215840906943
$ws.Range("K15").Value = 11.868299832869
$ws.Range("M15").Value = 14.12656765827226
$ws.Range("N15").Value = 17.23314332946067
$ws.Range("O15").Value = 17.66484923472413

$ws.Range("B16").Value = 6.610522752088456
$ws.Range("D16").Value = 4.304690003332726
$ws.Range("E16").Value = 13.68206358260758
$ws.Range("F16").Value = 20.07063420475066
$ws.Range("G16").Value = 21.69743290707055
$ws.Range("H16").Value = 12.20766004781364
$ws.Range("I16").Value = 18.3109313949038
$ws.Range("K16").Value = 11.64315987413496
$ws.Range("M16").Value = 13.94112526042299
$ws.Range("N16").Value = 17.26380198967747
$ws.Range("O16").Value = 17.67182326969386

$ws.Range("B17").Value = 6.556351331532463
$ws.Range("D17").Value = 4.270529415192512
$ws.Range("E17").Value = 13.60271017992
$ws.Range("F17").Value = 20.04279420830203
$ws.Range("G17").Value = 21.64309844019376
$ws.Range("H17").Value = 12.21809187411242
$ws.Range("I17").Value = 18.34776192541182
$ws.Range("K17").Value = 11.50238932588673
$ws.Range("M17").Value = 13.82681550294742
$ws.Range("N17").Value = 17.28303121497738
$ws.Range("O17").Value = 17.67725102039345

$ws.Range("B18").Value = 6.525064311673741
$ws.Range("D18").Value = 4.250692932886755
$ws.Range("E18").Value = 13.5573219826902
$ws.Range("F18").Value = 20.02733240245447
$ws.Range("G18").Value = 21.61258248118488
$ws.Range("H18").Value = 12.22429595550497
$ws.Range("I18").Value = 18.36923129572592
$ws.Range("K18").Value = 11.42046175415763
$ws.Range("M18").Value = 13.76089091101306
$ws.Range("N18").Value = 17.29424638359055
$ws.Range("O18").Value = 17.68079510496964

$ws.Range("B19").Value = 6.514449970365447
$ws.Range("D19").Value = 4.24394467000702
$ws.Range("E19").Value = 13.54199955067889
$ws.Range("F19").Value = 20.02219229120325
$ws.Range("G19").Value = 21.60237752309707
$ws.Range("H19").Value = 12.22643157359157
$ws.Range("I19").Value = 18.37654952852163
$ws.Range("K19").Value = 11.39255912634741
$ws.Range("M19").Value = 13.73854198896587
$ws.Range("N19").Value = 17.29807029809395
$ws.Range("O19").Value = 17.68206755709991

$ws.Range("B20").Value = 6.562131581221633
$ws.Range("D20").Value = 4.274185431574833
$ws.Range("E20").Value = 13.61113161082307
$ws.Range("F20").Value = 20.04570087915067
$ws.Range("G20").Value = 21.64880648920316
$ws.Range("H20").Value = 12.21696027594813
$ws.Range("I20").Value = 18.34381172447459
$ws.Range("K20").Value = 11.5174743154347
$ws.Range("M20").Value = 13.83900280971499
$ws.Range("N20").Value = 17.28096819237582
$ws.Range("O20").Value = 17.67662952675276

$ws.Range("B21").Value = 6.720716077463125
$ws.Range("D21").Value = 4.373516017802908
$ws.Range("E21").Value = 13.84631621632592
$ws.Range("F21").Value = 20.13124269158894
$ws.Range("G21").Value = 21.81352719919561
$ws.Range("H21").Value = 12.18757568592651
$ws.Range("I21").Value = 18.23722360686234
$ws.Range("K21").Value = 11.92559391658606
$ws.Range("M21").Value = 14.17425674723366
$ws.Range("N21").Value = 17.22535785008619
$ws.Range("O21").Value = 17.66340416462732

$ws.Range("B22").Value = 6.822981030094064
$ws.Range("D22").Value = 4.436682848231409
$ws.Range("E22").Value = 14.00191441540121
$ws.Range("F22").Value = 20.19197049497546
$ws.Range("G22").Value = 21.92753103338575
$ws.Range("H22").Value = 12.1701894362718
$ws.Range("I22").Value = 18.17009240851699
$ws.Range("K22").Value = 12.18346059293632
$ws.Range("M22").Value = 14.39128441446409
$ws.Range("N22").Value = 17.19038822568629
$ws.Range("O22").Value = 17.6585238344537

$ws.Range("B23").Value = 6.768537229233639
$ws.Range("D23").Value = 4.403133265224905
$ws.Range("E23").Value = 13.91871345084906
$ws.Range("F23").Value = 20.15912217780161
$ws.Range("G23").Value = 21.8661143751274
$ws.Range("H23").Value = 12.1793023455124
$ws.Range("I23").Value = 18.20569061600768
$ws.Range("K23").Value = 12.04665731156636
$ws.Range("M23").Value = 14.27566675384045
$ws.Range("N23").Value = 17.20892662985099
$ws.Range("O23").Value = 17.66078412154782

$ws.Range("B24").Value = 6.559518773608561
$ws.Range("D24").Value = 4.272533159648782
$ws.Range("E24").Value = 13.60732355162727
$ws.Range("F24").Value = 20.04438507851065
$ws.Range("G24").Value = 21.64622363016208
$ws.Range("H24").Value = 12.21747122768063
$ws.Range("I24").Value = 18.34559669126741
$ws.Range("K24").Value = 11.51065749258767
$ws.Range("M24").Value = 13.83349356541993
$ws.Range("N24").Value = 17.28190038634087
$ws.Range("O24").Value = 17.67690918458805

$ws.Range("B25").Value = 6.329279448572592
$ws.Range("D25").Value = 4.124516637681316
$ws.Range("E25").Value = 13.28135699106523
$ws.Range("F25").Value = 19.94172566418271
$ws.Range("G25").Value = 21.43725967092303
$ws.Range("H25").Value = 12.26643926932523
$ws.Range("I25").Value = 18.50742764926797
$ws.Range("K25").Value = 10.89603508756422
$ws.Range("M25").Value = 13.35012112415759
$ws.Range("N25").Value = 17.36654500929017
$ws.Range("O25").Value = 17.69699140797045
